$d = $word.ActiveDocument

# Reword degree lines to full, grammatically cleaner phrasing.
$d.Content.Find.Execute("MSc - Game Development", $true, $false, $false, $false, $false, $true, 1, $false, "Master’s degree in Game Development", 2)
$d.Content.Find.Execute("BSc - Software Engineering and Management", $true, $false, $false, $false, $false, $true, 1, $false, "Bachelor’s degree in Software Engineering and Management", 2)

# Fix "Eligible" -> "Elective" typo for every course table cell (replace all occurrences).
$d.Content.Find.Execute("(Eligible course)", $true, $false, $false, $false, $false, $true, 1, $false, "(Elective course)", 2)
